$p = $ppt.ActivePresentation

# --- Slide 10: "Independence of type names" -> two runs:
#     "Independent " + "of type names"
$s10 = $p.Slides.Item(10)
$sh10 = $s10.Shapes.Item(2)
$tr10 = $sh10.TextFrame.TextRange
$para10 = $tr10.Paragraphs(10, 1)
$run10_1 = $para10.Runs(1, 1)
$run10_1.Text = "Independent "
$run10_1.InsertAfter("of type names") | Out-Null

# --- Slide 2: "Problem " + "Development" (two runs) -> single run
#     "Problem Development"
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(2)
$tr2 = $sh2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(2, 1)
$run2_1 = $para2.Runs(1, 1)
$run2_1.Text = "Problem Development"
$run2_2 = $para2.Runs(2, 1)
$run2_2.Text = ""
